# "Organizando Localização das funções"
# Adds two new worksheets ("configs" and "VariableNames") to the workbook,
# populates them, and makes "VariableNames" the active sheet/tab.

$wb = $excel.ActiveWorkbook
$wsParams = $wb.Worksheets.Item(1)

# --- New sheet: configs (placed right after "params") ---------------------
$wsConfigs = $wb.Worksheets.Add($null, $wsParams)
$wsConfigs.Name = "configs"

$wsConfigs.Range("A1").Value = "Start"
$wsConfigs.Range("B1").Value = "Finish"
$wsConfigs.Range("C1").Value = "Step"
$wsConfigs.Range("A1").Font.Underline = $true

$wsConfigs.Range("A2").Value = 2015
$wsConfigs.Range("B2").Value = 2025
$wsConfigs.Range("C2").Value = 0.125

# --- New sheet: VariableNames (placed right after "configs") --------------
$wsVarNames = $wb.Worksheets.Add($null, $wsConfigs)
$wsVarNames.Name = "VariableNames"

$wsVarNames.Range("A1").Value = "ModelName"
$wsVarNames.Range("B1").Value = "ResultName"
$wsVarNames.Range("C1").Value = "ChartName"
$wsVarNames.Range("D1").Value = "TextName"
$wsVarNames.Range("E1").Value = "EquationName"

$wsVarNames.Range("B2").Value = "Tempo"
$wsVarNames.Range("C2").Value = "Anos"
$wsVarNames.Range("D2").Value = "Anos"

$wsVarNames.Range("B3").Value = "PotentialAdopters"
$wsVarNames.Range("C3").Value = "Clientes Potenciais"
$wsVarNames.Range("D3").Value = "Clientes Potenciais"

$wsVarNames.Range("B4").Value = "Adopters"
$wsVarNames.Range("C4").Value = "Clientes"
$wsVarNames.Range("D4").Value = "Clientes"

$wsVarNames.Range("B5").Value = "AdvEffectiveness"
$wsVarNames.Range("C5").Value = "Efetividade do Anúncio"
$wsVarNames.Range("D5").Value = "Efetividade do Anúncio"

$wsVarNames.Range("B6").Value = "ContactRate"
$wsVarNames.Range("C6").Value = "Taxa de Contato"
$wsVarNames.Range("D6").Value = "Taxa de Contato"

$wsVarNames.Range("B7").Value = "AdoptionFraction"
$wsVarNames.Range("C7").Value = "Taxa de Adoção"
$wsVarNames.Range("D7").Value = "Taxa de Adoção"

$wsVarNames.Range("B8").Value = "TotalPopulation"
$wsVarNames.Range("C8").Value = "População Total"
$wsVarNames.Range("D8").Value = "População Total"

$wsVarNames.Range("B9").Value = "Adoption_From_Advertising"
$wsVarNames.Range("C9").Value = "Novos Clientes em Propaganda"
$wsVarNames.Range("D9").Value = "Novos Clientes em Propaganda"

$wsVarNames.Range("B10").Value = "Adoption_From_Word_of_Mouth"
$wsVarNames.Range("C10").Value = "Novos Clientes por Boca a Boca"
$wsVarNames.Range("D10").Value = "Novos Clientes por Boca a Boca"

$wsVarNames.Range("B11").Value = "Adoption_Rate"
$wsVarNames.Range("C11").Value = "Taxa de novos Clientes"
$wsVarNames.Range("D11").Value = "Taxa de novos Clientes"

$wsVarNames.Range("B12").Value = "Replicacao"
$wsVarNames.Range("C12").Value = "Replicação"
$wsVarNames.Range("D12").Value = "Replicação"

$wsVarNames.Columns.Item(1).ColumnWidth = 11.2060546875
$wsVarNames.Columns.Item(2).ColumnWidth = 30.380859375
$wsVarNames.Columns.Item(3).ColumnWidth = 28.1396484375
$wsVarNames.Columns.Item(4).ColumnWidth = 28.1396484375
$wsVarNames.Columns.Item(5).ColumnWidth = 13.447265625

# --- Selection / active tab -------------------------------------------------
# "VariableNames" becomes the active sheet, with C12 selected; "params" loses
# its previous tabSelected/selection state (back to the sheet default).
$wsVarNames.Activate()
$wsVarNames.Range("C12").Select()
